# Metrics.xlsx update — "Extended Aged Care metrics - 7-day avg, % Weekly Change"
#
# The "Aged Care" metric block (previously rows 45-50, plus the
# "Treatments" block in rows 51-52) is expanded with new 7-day-average
# and Weekly-Change metrics. The net effect is that the table grows
# from 52 data rows to 59 data rows, i.e. ref A1:F52 -> A1:F59.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Target content for rows 45-59 (columns A: Metric Catergory,
# B: Metric Catergory - Sort, C: Metric, D: Metric - Sort).
# Rows flagged Style2 get the same "highlighted metric name" cell
# style that the existing "per 1M" rows already use.
# ---------------------------------------------------------------------
$rows = @(
    @{Row=45; A="Aged Care";  B=60; C="# Aged Care Resident Cases";                      D=450; Style2=$true},
    @{Row=46; A="Aged Care";  B=60; C="# Aged Care Resident Cases (7-day avg)";           D=460; Style2=$true},
    @{Row=47; A="Aged Care";  B=60; C="# Aged Care Resident Cases (7-day avg) per 1M";    D=470; Style2=$true},
    @{Row=48; A="Aged Care";  B=60; C="% Aged Care Resident Cases Weekly Change";         D=480; Style2=$true},
    @{Row=49; A="Aged Care";  B=60; C="# Aged Care Staff Cases";                          D=490; Style2=$true},
    @{Row=50; A="Aged Care";  B=60; C="# Aged Care Staff Cases (7-day avg)";              D=500; Style2=$true},
    @{Row=51; A="Aged Care";  B=60; C="# Aged Care Staff Cases (7-day avg) per 1M";       D=510; Style2=$true},
    @{Row=52; A="Aged Care";  B=60; C="% Aged Care Staff Cases Weekly Change";            D=520; Style2=$true},
    @{Row=53; A="Aged Care";  B=60; C="# Aged Care Active Outbreaks";                     D=530; Style2=$false},
    @{Row=54; A="Aged Care";  B=60; C="# Aged Care Active Outbreaks (7-day avg)";         D=540; Style2=$false},
    @{Row=55; A="Aged Care";  B=60; C="# Aged Care Active Outbreaks (7-day avg) per 1M";  D=550; Style2=$false},
    @{Row=56; A="Aged Care";  B=60; C="% Aged Care Active Outbreaks Weekly Change";       D=560; Style2=$false},
    @{Row=57; A="Treatments"; B=70; C="# Weekly PBS Scripts";                             D=570; Style2=$false},
    @{Row=58; A="Treatments"; B=70; C="# Weekly PBS Scripts per 1M";                      D=580; Style2=$false},
    @{Row=59; A="Treatments"; B=70; C="% Weekly PBS Scripts Change";                      D=590; Style2=$false}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D

    # column E ("covidlive.com.au") is always blank/centered here
    $ws.Cells.Item($r.Row, 5).Value = ""
    $ws.Cells.Item($r.Row, 5).HorizontalAlignment = -4108

    # column F ("health.gov.au") carries the "X" marker, centered
    $ws.Cells.Item($r.Row, 6).Value = "X"
    $ws.Cells.Item($r.Row, 6).HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------
# Apply the highlighted-metric style (matching the existing "per 1M"
# rows, e.g. C47) to every Style2-flagged row's Metric cell.
# ---------------------------------------------------------------------
$styleSource = $ws.Range("C47")
$styleSource.Copy()
foreach ($r in $rows) {
    if ($r.Style2) {
        $ws.Cells.Item($r.Row, 3).PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Grow the "Metrics" table / autofilter to cover the new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F59"))

# Column C widened to fit the longer metric names now present.
$ws.Columns.Item(3).ColumnWidth = 41

# ---------------------------------------------------------------------
# Refresh the on-screen view to roughly match where the edit happened.
# ---------------------------------------------------------------------
$null = $ws.Range("D43:D59").Select()
$excel.ActiveWindow.ScrollRow = 35
